# Refresh the per-coin Price (D) / Volume(1h) (E) figures pulled by the
# GitHub Actions scraper. Values are stored as literal text in the sheet
# (matching the original inline-string cells), so numeric-looking Price
# values are entered with a leading apostrophe to force Excel to keep
# them as text instead of auto-converting to a number (which would drop
# trailing zeros / exact formatting, e.g. "151.50" -> 151.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.336.51"
$ws.Range("E2").Value = "  -3.68%  "
$ws.Range("D3").Value = "1.852.29"
$ws.Range("E3").Value = "  -5.13%  "
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'322.66"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.4496"
$ws.Range("E7").Value = "  -5.57%  "
$ws.Range("D8").Value = "'0.3832"
$ws.Range("E8").Value = "  -4.74%  "
$ws.Range("D9").Value = "'48.33"
$ws.Range("E9").Value = "  -9.76%  "
$ws.Range("D10").Value = "'0.07853"
$ws.Range("E10").Value = "  -6.78%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D12").Value = "'21.32"
$ws.Range("E12").Value = "  -4.20%  "
$ws.Range("D13").Value = "1.814.59"
$ws.Range("E13").Value = "  -7.33%  "
$ws.Range("D14").Value = "'5.855"
$ws.Range("E14").Value = "  -4.60%  "
$ws.Range("D15").Value = "'7.127"
$ws.Range("E15").Value = "  -5.53%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'85.71"
$ws.Range("E17").Value = "  -5.37%  "
$ws.Range("D18").Value = "'0.00001030"
$ws.Range("D19").Value = "'0.06498"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").Value = "'16.98"
$ws.Range("E20").Value = "  -7.75%  "
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("E22").Value = "  -5.91%  "
$ws.Range("D23").Value = "27.315.24"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("E24").Value = "  -5.79%  "
$ws.Range("D25").Value = "'2.263"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "2.053.87"
$ws.Range("E26").Value = "  -6.17%  "
$ws.Range("D27").Value = "'151.50"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("E28").Value = "  -3.95%  "
$ws.Range("D29").Value = "'5.542"
$ws.Range("E29").Value = "  -6.23%  "
$ws.Range("D30").Value = "'2.051"
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("D31").Value = "'119.84"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("D32").Value = "'0.09320"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("D34").Value = "'0.9336"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").Value = "'3.594"
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").Value = "'5.258"
$ws.Range("E36").Value = "  -5.93%  "
$ws.Range("D37").Value = "'0.02222"
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("D38").Value = "'0.05977"
$ws.Range("E38").Value = "  -3.79%  "
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("D40").Value = "'8.276"
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").Value = "'0.5893"
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("E44").Value = "  -8.11%  "
$ws.Range("D45").Value = "'1.251"
$ws.Range("E45").Value = "  -6.54%  "
$ws.Range("D46").Value = "'0.5640"
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("D47").Value = "'12.18"
$ws.Range("E47").Value = "  -6.12%  "
$ws.Range("E48").Value = "  -6.17%  "
$ws.Range("D49").Value = "'3.360"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'0.06874"
$ws.Range("E50").Value = "  +1.15%  "
